$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hike Difficulties")
$lo = $ws.ListObjects.Item(1)

$lr1 = $lo.ListRows.Add()
$lr1.Range.Cells.Item(1,1).Value2 = "Deadwood Lakes"

$lr2 = $lo.ListRows.Add()
$lr2.Range.Cells.Item(1,1).Value2 = "Silver Forest Trail"

$lr3 = $lo.ListRows.Add()
$lr3.Range.Cells.Item(1,1).Value2 = "Naches Peak Loop"

$lr1.Range.Cells.Item(1,2).Value2 = 3.4
$lr1.Range.Cells.Item(1,3).Value2 = 940
$lr1.Range.Cells.Item(1,4).Value2 = "moderate (because rough)"

$lr2.Range.Cells.Item(1,2).Value2 = 3.3
$lr2.Range.Cells.Item(1,3).Value2 = 600
$lr2.Range.Cells.Item(1,4).Value2 = "easy"

$lr3.Range.Cells.Item(1,2).Value2 = 4.2
$lr3.Range.Cells.Item(1,3).Value2 = 860
$lr3.Range.Cells.Item(1,4).Value2 = "easy"

$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($lo.ListColumns.Item(1).Range)
$lo.Sort.Header = 1
$lo.Sort.Apply()

Write-Host $lo.Range.Address()
